$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the SECOND occurrence of "7947" in the document (the first
# lives in the SVMs results table and must stay untouched; the second
# lives in the "Final F1 Score" summary table at the very end of the
# document and is the one edited in this revision).
# ------------------------------------------------------------------
$probe = $d.Content
$found1 = $probe.Find.Execute("7947", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$probe.Collapse(0)
$found2 = $probe.Find.Execute("7947", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitStart = $probe.Start
$splitEnd = $probe.End

# ------------------------------------------------------------------
# Force a run boundary right before the text we are about to retype
# (mirrors the original author selecting "7947" and typing "8134" -
# the untouched "0." prefix stays in its own run while the newly
# typed "8134" becomes a second run).
# ------------------------------------------------------------------
$splitPointRange = $d.Range($splitStart, $splitStart)
$d.Bookmarks.Add("ZZZ_SPLIT_TMP", $splitPointRange)

# Replace "7947" with "8134"
$editRange = $d.Range($splitStart, $splitEnd)
$editRange.Text = "8134"

# ------------------------------------------------------------------
# Relocate the "_GoBack" bookmark (Word always keeps exactly one,
# tracking the location of the most recent edit) to the end of the
# text we just typed. Adding a bookmark named "_GoBack" replaces any
# existing one elsewhere in the document.
# ------------------------------------------------------------------
$newGoBackPoint = $d.Range($splitStart + 4, $splitStart + 4)
$d.Bookmarks.Add("_GoBack", $newGoBackPoint)

# Drop the temporary split marker now that its job (forcing the run
# break) is done - only "_GoBack" should remain.
$d.Bookmarks("ZZZ_SPLIT_TMP").Delete()
